$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.865433692365116
$ws.Range("C2").Value = 0.1129504528410621
$ws.Range("D2").Value = 0.2588888633272575
$ws.Range("E2").Value = 0.05888425167063538
$ws.Range("F2").Value = 4.166620448711939
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("L2").Value = 0.1950833067680051
$ws.Range("M2").Value = 0.3646298907252614
$ws.Range("B3").Value = 1.801585411290148
$ws.Range("C3").Value = 0.0977230480798994
$ws.Range("D3").Value = 0.2485202626132548
$ws.Range("E3").Value = 0.05833157205361506
$ws.Range("F3").Value = 3.966901955159841
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("L3").Value = 0.1928430019437073
$ws.Range("M3").Value = 0.3549506182449633
$ws.Range("B4").Value = 1.764192650805796
$ws.Range("C4").Value = 0.08840518672499798
$ws.Range("D4").Value = 0.242147989246277
$ws.Range("E4").Value = 0.05797938699968741
$ws.Range("F4").Value = 3.845285543863554
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("L4").Value = 0.1915707671165308
$ws.Range("M4").Value = 0.3493223301951289
$ws.Range("B5").Value = 1.749407672288896
$ws.Range("C5").Value = 0.08461533282630285
$ws.Range("D5").Value = 0.2395492152180196
$ws.Range("E5").Value = 0.05783260064830364
$ws.Range("F5").Value = 3.795972201350452
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("L5").Value = 0.1910782396036126
$ws.Range("M5").Value = 0.3471075890278783
$ws.Range("B6").Value = 1.746979923001163
$ws.Range("C6").Value = 0.083986446589563
$ws.Range("D6").Value = 0.2391175531486738
$ws.Range("E6").Value = 0.05780802833911958
$ws.Range("F6").Value = 3.787798403158746
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("L6").Value = 0.1909980194219827
$ws.Range("M6").Value = 0.3467445846585164
$ws.Range("B7").Value = 1.763991424786809
$ws.Range("C7").Value = 0.08835404706127292
$ws.Range("D7").Value = 0.2421129501593242
$ws.Range("E7").Value = 0.05797742067276257
$ws.Range("F7").Value = 3.844619500044558
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("L7").Value = 0.19156401985542
$ws.Range("M7").Value = 0.349292142664595
$ws.Range("B8").Value = 1.843041607641396
$ws.Range("C8").Value = 0.1076929774053212
$ws.Range("D8").Value = 0.2553146713279233
$ws.Range("E8").Value = 0.05869632796336399
$ws.Range("F8").Value = 4.097542649461417
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("L8").Value = 0.1942893658410725
$ws.Range("M8").Value = 0.3612269279132434
$ws.Range("B9").Value = 2.012543560457118
$ws.Range("C9").Value = 0.1459051811705478
$ws.Range("D9").Value = 0.281181525930208
$ws.Range("E9").Value = 0.0600060731653369
$ws.Range("F9").Value = 4.601939369511427
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("L9").Value = 0.2004570099455378
$ws.Range("M9").Value = 0.3871463492379164
$ws.Range("B10").Value = 2.14609497846368
$ws.Range("C10").Value = 0.1742087125243756
$ws.Range("D10").Value = 0.3002092042610514
$ws.Range("E10").Value = 0.06090993594626504
$ws.Range("F10").Value = 4.978227210593531
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("L10").Value = 0.2054959191073493
$ws.Range("M10").Value = 0.4077501761971263
$ws.Range("B11").Value = 2.208852338371457
$ws.Range("C11").Value = 0.1871463113645859
$ws.Range("D11").Value = 0.3088784807000877
$ws.Range("E11").Value = 0.06130900693438357
$ws.Range("F11").Value = 5.15077811869844
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("L11").Value = 0.2078997798860911
$ws.Range("M11").Value = 0.4174686230511355
$ws.Range("B12").Value = 2.232908395317736
$ws.Range("C12").Value = 0.1920553194998718
$ws.Range("D12").Value = 0.3121638770571735
$ws.Range("E12").Value = 0.06145843032003051
$ws.Range("F12").Value = 5.216326338398346
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("L12").Value = 0.2088262092430284
$ws.Range("M12").Value = 0.4211989285676196
$ws.Range("B13").Value = 2.227714492553332
$ws.Range("C13").Value = 0.1909976258157826
$ws.Range("D13").Value = 0.3114561862072662
$ws.Range("E13").Value = 0.06142632405104331
$ws.Range("F13").Value = 5.202200006244084
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("L13").Value = 0.2086259668894002
$ws.Range("M13").Value = 0.4203933034372511
$ws.Range("B14").Value = 2.210825590846923
$ws.Range("C14").Value = 0.187549976504215
$ws.Range("D14").Value = 0.3091487181374646
$ws.Range("E14").Value = 0.06132133382548788
$ws.Range("F14").Value = 5.156166610150819
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("L14").Value = 0.2079756738222329
$ws.Range("M14").Value = 0.4177745100974448
$ws.Range("B15").Value = 2.200518665624259
$ws.Range("C15").Value = 0.1854394970397379
$ws.Range("D15").Value = 0.3077356743565076
$ws.Range("E15").Value = 0.06125680470785833
$ws.Range("F15").Value = 5.127997048010911
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("L15").Value = 0.2075794549929668
$ws.Range("M15").Value = 0.4161769664950867
$ws.Range("B16").Value = 2.14203422139343
$ws.Range("C16").Value = 0.1733645359396121
$ws.Range("D16").Value = 0.2996429683624626
$ws.Range("E16").Value = 0.06088361596336034
$ws.Range("F16").Value = 4.966979060952013
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("L16").Value = 0.2053410748586231
$ws.Range("M16").Value = 0.4071220452936259
$ws.Range("B17").Value = 2.106671333836118
$ws.Range("C17").Value = 0.1659734837244002
$ws.Range("D17").Value = 0.2946822278969989
$ws.Range("E17").Value = 0.06065160854612772
$ws.Range("F17").Value = 4.868558525207504
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("L17").Value = 0.2039965541798523
$ws.Range("M17").Value = 0.4016559817493359
$ws.Range("B18").Value = 2.086519990484135
$ws.Range("C18").Value = 0.161728145016923
$ws.Range("D18").Value = 0.2918301752711159
$ws.Range("E18").Value = 0.06051702237566947
$ws.Range("F18").Value = 4.812078568048264
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("L18").Value = 0.203233723771973
$ws.Range("M18").Value = 0.3985445551131619
$ws.Range("B19").Value = 2.079729368793608
$ws.Range("C19").Value = 0.1602917157978254
$ws.Range("D19").Value = 0.2908647133627511
$ws.Range("E19").Value = 0.06047125640477669
$ws.Range("F19").Value = 4.792977293862236
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("L19").Value = 0.2029772437074939
$ws.Range("M19").Value = 0.3974966498716412
$ws.Range("B20").Value = 2.11041624776766
$ws.Range("C20").Value = 0.1667596673647154
$ws.Range("D20").Value = 0.2952101752153737
$ws.Range("E20").Value = 0.06067642397075135
$ws.Range("F20").Value = 4.879022149741331
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("L20").Value = 0.2041385931356956
$ws.Range("M20").Value = 0.4022344859311247
$ws.Range("B21").Value = 2.215778342925432
$ws.Range("C21").Value = 0.1885623610651521
$ws.Range("D21").Value = 0.3098264040602601
$ws.Range("E21").Value = 0.06135221767937948
$ws.Range("F21").Value = 5.169682044903141
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("L21").Value = 0.2081662419351318
$ws.Range("M21").Value = 0.4185423488439923
$ws.Range("B22").Value = 2.286337404458834
$ws.Range("C22").Value = 0.2028693429196267
$ws.Range("D22").Value = 0.3193939772040153
$ws.Range("E22").Value = 0.06178402237172986
$ws.Range("F22").Value = 5.360855604433425
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("L22").Value = 0.2108926587162756
$ws.Range("M22").Value = 0.4294929498645885
$ws.Range("B23").Value = 2.248522241251578
$ws.Range("C23").Value = 0.1952278662002698
$ws.Range("D23").Value = 0.3142860121293438
$ws.Range("E23").Value = 0.06155444810972432
$ws.Range("F23").Value = 5.258708865385358
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("L23").Value = 0.2094288778196756
$ws.Range("M23").Value = 0.4236214996045149
$ws.Range("B24").Value = 2.108722614076385
$ws.Range("C24").Value = 0.1664042217356325
$ws.Range("D24").Value = 0.2949714904537188
$ws.Range("E24").Value = 0.06066520866456271
$ws.Range("F24").Value = 4.874291223958579
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("L24").Value = 0.2040743457080367
$ws.Range("M24").Value = 0.4019728473969124
$ws.Range("B25").Value = 1.965118975825703
$ws.Range("C25").Value = 0.1355314440097288
$ws.Range("D25").Value = 0.27418242419607
$ws.Range("E25").Value = 0.05966219706069231
$ws.Range("F25").Value = 4.464526573961734
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("L25").Value = 0.1986997894115561
$ws.Range("M25").Value = 0.3798623891581911
